$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.661.15'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '3.497.20'
$ws.Range('E3').Value = '  -2.89%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.90%  '
$ws.Range('D7').Value = '3.493.33'
$ws.Range('E7').Value = '  -2.95%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.504'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.76'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.37%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.128'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -6.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.397'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.59%  '
$ws.Range('D13').Value = '4.087.63'
$ws.Range('E13').Value = '  -2.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000190'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -9.41%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -7.88%  '
$ws.Range('D16').Value = '3.494.24'
$ws.Range('E16').Value = '  -2.99%  '
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '65.499.48'
$ws.Range('E18').Value = '  -1.81%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.08'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.27'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '412.34'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.584'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.88%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '76.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.74%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.629.50'
$ws.Range('E26').Value = '  -2.82%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000110'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.61%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.42'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -7.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.67'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.16%  '
$ws.Range('E31').Value = '  -0.06%  '
$ws.Range('D32').Value = '3.494.14'
$ws.Range('E32').Value = '  -2.77%  '
$ws.Range('E33').Value = '  -2.26%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.83'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.32'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '175.51'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('E39').Value = '  -8.93%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.09'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -9.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0799'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -6.98%  '
$ws.Range('B42').Value = 'Mantle'
$ws.Range('C42').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.846'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.84%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.85'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '45.21'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.14%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.997'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.73'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -9.60%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '22.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.74%  '
$ws.Range('B48').Value = 'dogwifhat'
$ws.Range('C48').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.28'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -11.77%  '
$ws.Range('B49').Value = 'Cosmos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.44%  '
$ws.Range('E50').Value = '  -9.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.885'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.12%  '
